$d = $word.ActiveDocument

function Replace-Text($range, $old, $new) {
    $found = $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Output "WARNING: not found: $old"
    }
    return $found
}

# ------------------------------------------------------------------
# Title / byline / email (scope email-domain replace to its own
# paragraph so the short word "com" cannot clobber other occurrences)
# ------------------------------------------------------------------
[void](Replace-Text $d.Content "Cyber Ecosystem's Balance: Harmony in the Digital Realm" "The Allure of Government: A Path Forward")
[void](Replace-Text $d.Content "Katelyn Scholes" "Anais Cooper")
[void](Replace-Text $d.Content "katelynscholes@ritemail" "anaiscooper@edumail")
[void](Replace-Text $d.Paragraphs.Item(3).Range "com" "org")

# ------------------------------------------------------------------
# First body paragraph, section 1 (before the first double line-break)
# ------------------------------------------------------------------
[void](Replace-Text $d.Content "In the interconnected expanse of cyberspace, where information flows like an ever-present river, a delicate balance reigns, shaping the digital landscape" "Government, an intricate tapestry woven by human hands, shapes our destiny and weaves the fabric of our societies")

[void](Replace-Text $d.Content " This equilibrium, intricately woven into the fabric of technology, is a testament to the harmonious coexistence of diverse components, each playing a pivotal role in maintaining the stability of the virtual world" " From the dawn of civilization, humans have yearned for systems that can regulate their interactions, foster cooperation, and uphold justice")

$rng = $d.Content
$ok = Replace-Text $rng " From the complex interactions between firewalls and malware to the interplay of regulations and user behavior, the cyber ecosystem thrives amidst constant fluctuations, adapting to the ebb and flow of innovation and challenges" " Government, in its myriad forms, embodies this desire for order and progress"
if ($ok) {
    $rng.Collapse(0)
    $rng.InsertAfter(".")
    $rng.Collapse(0)
    $rng.InsertAfter(" Understanding government is not merely an academic exercise; it is a profound journey into the heart of human nature and the dynamics that shape our collective existence")
}

# ------------------------------------------------------------------
# First body paragraph, section 2 (between the two double line-breaks)
# ------------------------------------------------------------------
[void](Replace-Text $d.Content "As technology continues its relentless march forward, the boundaries of the cyber ecosystem expand, mirroring the rapid evolution of the real world" "Governments exist in a myriad of forms, each reflecting the unique history, culture, and aspirations of its people")

[void](Replace-Text $d.Content " New technologies emerge, reshaping the digital landscape, while old ones fade into obsolescence" " Democracy, with its emphasis on popular sovereignty and representative governance, stands as a beacon of inclusion and empowerment")

[void](Replace-Text $d.Content " Within this ever-shifting panorama, the balance remains a guiding principle, dictating the trajectories of progress and influencing the impact of digital transformation on society. Understanding and nurturing this equilibrium are paramount to ensuring a resilient and sustainable cyber ecosystem, one that can withstand the complexities and uncertainties of the digital age" " In democratic societies, citizens actively participate in the decision-making process, shaping policies and electing leaders entrusted with the responsibility of governing")

# ------------------------------------------------------------------
# First body paragraph, section 3 (after the second double line-break)
# ------------------------------------------------------------------
[void](Replace-Text $d.Content "The interplay between technology, regulation, and human behavior forms the cornerstone of the cyber ecosystem's balance" "Government is a dynamic and evolving entity, adapting to the ever-changing needs and aspirations of its citizens")

[void](Replace-Text $d.Content " Technological advancements, acting as catalysts for innovation, drive the development of new solutions that enhance user experiences and expand the possibilities of digital interaction" " As societies progress, governments must grapple with new challenges and opportunities")

[void](Replace-Text $d.Content " Regulations, crafted with the intention of safeguarding users and upholding democratic values, establish boundaries and provide a framework for responsible innovation. Human behavior, often unpredictable and multifaceted, introduces an element of dynamic unpredictability into the equation, shaping the actualization of technology and the effectiveness of regulation. Striking a balance between these three pillars is a continuous endeavor, requiring collaboration, adaptability, and a deep understanding of the evolving digital realm" " They must find ways to promote economic growth while ensuring social equity, protect the environment while fostering development, and navigate the complex intricacies of international relations")

# ------------------------------------------------------------------
# "Summary" heading: re-assert its own text to drop the stale
# lastRenderedPageBreak bookkeeping element, then the body text.
# ------------------------------------------------------------------
[void](Replace-Text $d.Content "Summary" "Summary")

[void](Replace-Text $d.Content "The cyber ecosystem is a complex interconnected system where information flows like an ever-present river" "In conclusion, government is a multifaceted entity that plays a crucial role in shaping our lives and societies")

[void](Replace-Text $d.Content " The delicate balance of this ecosystem is shaped by the interplay of technology, regulation, and human behavior" " By understanding the diverse forms of government, the principles of democratic governance, and the challenges and opportunities that governments face, we gain insights into the intricate tapestry of human civilization")

[void](Replace-Text $d.Content " Understanding and nurturing this equilibrium is paramount to ensuring a resilient and sustainable cyber ecosystem, one that can withstand the complexities and uncertainties of the digital age" " The study of government is not just an academic pursuit but a profound exploration of human nature and the art of living together")

[void](Replace-Text $d.Content " Striking a balance between technological advancements, regulatory frameworks, and human behavior is a continuous endeavor, requiring collaboration, adaptability, and a deep understanding of the evolving digital realm" " As we navigate the complexities of modern governance, may we strive to create governments that uphold justice, promote equality, and empower citizens to shape their destiny")

# ------------------------------------------------------------------
# Append a new empty paragraph at the very end of the document
# ------------------------------------------------------------------
$endRng = $d.Content
$endRng.Collapse(0)
$endRng.InsertParagraphAfter()

Write-Output "DONE"
